$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the Air_UAV rows (211-212) down to 216-217 first, before overwriting ---
$ws.Range("A216").Value = $ws.Range("A211").Value2
$ws.Range("B216").Value = $ws.Range("B211").Value2
$ws.Range("C216").Value = $ws.Range("C211").Value2
$ws.Range("C216").HorizontalAlignment = $ws.Range("C211").HorizontalAlignment

$ws.Range("A217").Value = $ws.Range("A212").Value2
$ws.Range("B217").Value = $ws.Range("B212").Value2
$ws.Range("C217").Value = $ws.Range("C212").Value2
$ws.Range("C217").HorizontalAlignment = $ws.Range("C212").HorizontalAlignment

# --- Row 205: early_helicopter moves from A to B; add note in D ---
$ws.Range("A205").ClearContents()
$ws.Range("B205").Value = "early_helicopter"
$ws.Range("D205").Value = "required for all helicopter research"

# --- Row 206-209: Attack Helicopter tech codes added in column B ---
$ws.Range("B206").Value = "attack_helicopter1"
$ws.Range("B207").Value = "attack_helicopter2"
$ws.Range("B208").Value = "attack_helicopter3"
$ws.Range("B209").Value = "attack_helicopter4"

# --- Rows 211-214: Transport Helicopter section (overwrites the old Air_UAV data, now moved) ---
$ws.Range("A211").Value = "transport_helicopter_equipment_1"
$ws.Range("B211").Value = "transport_helicopter1"
$ws.Range("C211").Value = 1965
$ws.Range("D211").Value = "Transport Helicopter"
$ws.Range("D211").HorizontalAlignment = 1

$ws.Range("A212").Value = "transport_helicopter_equipment_2"
$ws.Range("B212").Value = "transport_helicopter2"
$ws.Range("C212").Value = 1985

$ws.Range("A213").Value = "transport_helicopter_equipment_3"
$ws.Range("B213").Value = "transport_helicopter3"
$ws.Range("C213").Value = 2005

$ws.Range("A214").Value = "transport_helicopter_equipment_4"
$ws.Range("B214").Value = "transport_helicopter4"
$ws.Range("C214").Value = 2015

# The new Transport Helicopter block's Generation column (C) keeps the default/general
# alignment (style differs from the rest of the sheet's centered Generation column).
$ws.Range("C211:C214").HorizontalAlignment = 1
